# New submission synced: 2026-02-04 18:56:30
# Sheet "JSS 3C" (Project_Results.xlsx) gets its Admission-No column cell
# C2 corrected from a text "42" to a genuine number, and a brand new
# response row (row 3) appended for NAMBA MARCUS PULKA.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3C")

# --- fix existing row 2: Admission No should be numeric, not text ---
$ws.Range("C2").Value = 42

# --- append new submission as row 3 ---
$ws.Range("A3").Value = "2026-02-04 18:56:30"
$ws.Range("B3").Value = "NAMBA MARCUS PULKA"

# Admission No for this submission stays textual (matches form's raw
# string answer) - use a leading apostrophe to force text, then drop the
# resulting "quote prefix" style so the cell keeps the sheet's default
# formatting.
$ws.Range("C3").Value = "'39"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = 9
